$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.928.37"
$ws.Range("E2").Value = "  +1.39%  "
$ws.Range("D3").Value = "2.408.41"
$ws.Range("E3").Value = "  +1.62%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "'555.43"
$ws.Range("E5").Value = "  +1.13%  "
$ws.Range("D6").Value = "'142.01"
$ws.Range("E6").Value = "  +2.94%  "
$ws.Range("D8").Value = "'0.529"
$ws.Range("E8").Value = "  +0.46%  "
$ws.Range("D9").Value = "2.403.19"
$ws.Range("E9").Value = "  +1.37%  "
$ws.Range("E10").Value = "  +0.95%  "
$ws.Range("E11").Value = "  -1.12%  "
$ws.Range("E12").Value = "  +0.56%  "
$ws.Range("D13").Value = "'0.351"
$ws.Range("E13").Value = "  +0.98%  "
$ws.Range("D14").Value = "'25.89"
$ws.Range("E14").Value = "  +3.59%  "
$ws.Range("E15").Value = "  +4.02%  "
$ws.Range("D16").Value = "2.839.50"
$ws.Range("E16").Value = "  +2.15%  "
$ws.Range("D17").Value = "61.877.70"
$ws.Range("E17").Value = "  +1.47%  "
$ws.Range("D18").Value = "2.403.37"
$ws.Range("E18").Value = "  +0.85%  "
$ws.Range("D19").Value = "'11.07"
$ws.Range("E19").Value = "  +2.63%  "
$ws.Range("B20").Value = "Polkadot"
$ws.Range("C20").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D20").Value = "'4.16"
$ws.Range("E20").Value = "  +0.53%  "
$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D21").Value = "'322.21"
$ws.Range("E21").Value = "  +0.53%  "
$ws.Range("E22").Value = "  +0.36%  "
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("D24").Value = "'65.00"
$ws.Range("E24").Value = "  +1.33%  "
$ws.Range("E25").Value = "  +1.97%  "
$ws.Range("D26").Value = "'9.17"
$ws.Range("E26").Value = "  +10.06%  "
$ws.Range("D27").Value = "'577.16"
$ws.Range("E27").Value = "  +14.54%  "
$ws.Range("B28").Value = "WrappedeETH"
$ws.Range("C28").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D28").Value = "2.644.29"
$ws.Range("E28").Value = "  +6.83%  "
$ws.Range("B29").Value = "Binance-PegBSC-USD"
$ws.Range("C29").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "  +0.27%  "
$ws.Range("D30").Value = "'8.23"
$ws.Range("E30").Value = "  +1.20%  "
$ws.Range("D31").Value = "0.0₃0921"
$ws.Range("E31").Value = "  +4.15%  "
$ws.Range("D32").Value = "'1.44"
$ws.Range("E32").Value = "  +3.91%  "
$ws.Range("E33").Value = "  -1.52%  "
$ws.Range("D34").Value = "'1.87"
$ws.Range("E34").Value = "  +2.20%  "
$ws.Range("E35").Value = "  +2.30%  "
$ws.Range("D37").Value = "'5.62"
$ws.Range("E37").Value = "  +5.26%  "
$ws.Range("E38").Value = "  +0.56%  "
$ws.Range("E39").Value = "  +0.78%  "
$ws.Range("D40").Value = "'151.26"
$ws.Range("E40").Value = "  +3.75%  "
$ws.Range("D41").Value = "'18.62"
$ws.Range("E41").Value = "  +0.37%  "
$ws.Range("D42").Value = "'1.82"
$ws.Range("E42").Value = "  -3.59%  "
$ws.Range("D43").Value = "'0.999"
$ws.Range("E43").Value = "  +0.01%  "
$ws.Range("E44").Value = "  +11.65%  "
$ws.Range("D45").Value = "'149.71"
$ws.Range("E45").Value = "  +1.32%  "
$ws.Range("D46").Value = "'3.62"
$ws.Range("E46").Value = "  +0.74%  "
$ws.Range("D47").Value = "'0.0537"
$ws.Range("E47").Value = "  +3.25%  "
$ws.Range("D48").Value = "'20.07"
$ws.Range("E48").Value = "  +3.73%  "
$ws.Range("E49").Value = "  +2.11%  "
$ws.Range("D50").Value = "'0.0920"
$ws.Range("E50").Value = "  +1.24%  "
$ws.Range("E51").Value = "  +1.70%  "
